$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.93%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.71%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.554"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.68%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08481"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.77%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.080"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9898"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.59%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.29%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1148"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.23%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1938"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.42%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "9.484"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.85%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09834"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.35%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04686"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.74%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1060"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.15%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001295"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.24%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005883"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.22%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.004626"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "6.10%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.386"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.22%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.430"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.77%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3354"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.48%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1383"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.05%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2552"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.02%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04143"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.30%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001302"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.02%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.29%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02727"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.60%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05754"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007820"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.19%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1434"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.29%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007253"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.19%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.94%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008055"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3555"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007060"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.69%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.19%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.27%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003538"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.11%"
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003425"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-2.23%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.19%"
